$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94 (shifts existing rows 94:123 down to 95:124)
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new data record
$ws.Range("A94").Value = 10
$ws.Range("B94").Value = "Vega Modelo de Temuco"
$ws.Range("C94").Value = "La Araucanía"
$ws.Range("D94").Value = 45141
$ws.Range("E94").Value = 9
$ws.Range("F94").Value = 300000001
$ws.Range("G94").Value = "Rabanito"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 40
$ws.Range("K94").Value = 8000
$ws.Range("L94").Value = 8000
$ws.Range("M94").Value = 8000
$ws.Range("N94").Value = "$/docena de paquetes"
$ws.Range("O94").Value = "Provincia de Cautín"
$ws.Range("P94").Value = 667
$ws.Range("Q94").Value = 12
$ws.Range("R94").Value = "Hortaliza"
